# Insert a new weekly price record for "Femacal de La Calera - Berenjena"
# at row 305, pushing the existing rows 305:371 down to 306:372.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 305..371 down by one row.
$ws.Rows("305").Insert()

# Populate the newly inserted row 305 with the new record.
$ws.Range("A305").Value = 3
$ws.Range("B305").Value = "Femacal de La Calera"
$ws.Range("C305").Value = "Coquimbo"
$ws.Range("D305").Value = 44855
$ws.Range("E305").Value = 5
$ws.Range("F305").Value = 100112001
$ws.Range("G305").Value = "Berenjena"
$ws.Range("H305").Value = "Sin especificar"
$ws.Range("I305").Value = "Primera"
$ws.Range("J305").Value = 105
$ws.Range("K305").Value = 11500
$ws.Range("L305").Value = 12000
$ws.Range("M305").Value = 11762
$ws.Range("N305").Value = "$/caja 60 unidades"
$ws.Range("O305").Value = "Región de Arica y Parinacota"
$ws.Range("P305").Value = 196
$ws.Range("Q305").Value = 60
$ws.Range("R305").Value = "Hortaliza"

# Match the date-format style used by the other rows' Fecha (D) column.
$ws.Range("D305").NumberFormat = $ws.Range("D306").NumberFormat
